# Updated symbol list on Sat Dec 17 02:36:16 UTC 2022 with GitHub Actions
#
# Refreshes the crypto price/volume table on Sheet1. Every touched cell in
# column D holds a price that is stored as TEXT (not a number) in the
# workbook, so we force the "@" text format before writing the new value -
# this keeps things like trailing zeros ("5.100") and long decimals
# ("0.00005445") intact instead of having Excel silently normalise them to
# a float. Columns B/C/E are plain strings already and need no such care.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $value) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $value
}

# Row 2 - BNB
Set-TextValue "D2" "224.22"

# Row 3 - OKB
Set-TextValue "D3" "22.46"

# Row 4 - HuobiToken
Set-TextValue "D4" "5.100"

# Row 5 - Cronos
Set-TextValue "D5" "0.05528"

# Row 6 - GateToken
Set-TextValue "D6" "3.391"

# Row 8 - was FTXToken, now MXToken (swapped with row 9)
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue "D8" "0.7827"
$ws.Range("E8").Value = "7MXTokenMX"

# Row 9 - was MXToken, now FTXToken (swapped with row 8)
$ws.Range("B9").Value = "FTXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
Set-TextValue "D9" "1.036"
$ws.Range("E9").Value = "8FTXTokenFTTWorstin24h"

# Row 10 - WazirX
Set-TextValue "D10" "0.1374"

# Row 11 - MandalaExchangeToken
Set-TextValue "D11" "0.07338"

# Row 12 - LiechtensteinCryptoassetsExchange
Set-TextValue "D12" "0.03124"

# Row 14 - BitMartToken
Set-TextValue "D14" "0.09243"

# Row 15 - BitForexToken
Set-TextValue "D15" "0.001665"

# Row 16 - MCDex
Set-TextValue "D16" "3.260"

# Row 17 - CoinExToken
Set-TextValue "D17" "0.04772"

# Row 18 - One
Set-TextValue "D18" "0.0005874"

# Row 19 - TigerCash
Set-TextValue "D19" "0.006263"

# Row 20 - HotbitToken
Set-TextValue "D20" "0.005242"

# Row 21 - BitKan
Set-TextValue "D21" "0.001065"

# Row 22 - NitroEx
Set-TextValue "D22" "0.0001502"

# Row 23 - LEO
Set-TextValue "D23" "3.838"

# Row 24 - BTSEToken
Set-TextValue "D24" "2.196"

# Row 26 - ProBitToken
Set-TextValue "D26" "0.1288"

# Row 27 - UpBots (volume label gained "Bestin24h" suffix)
$ws.Range("E27").Value = "26UpBotsUBXTBestin24h"

# Row 40 - IDEX
Set-TextValue "D40" "0.03872"

# Row 41 - KickToken (volume label lost "Bestin24h" suffix)
Set-TextValue "D41" "0.007126"
$ws.Range("E41").Value = "40KickTokenKICK"

# Row 42 - was CEJI, now BKEXToken (swapped with row 43)
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D42" "0.1030"
$ws.Range("E42").Value = "41BKEXTokenBKK"

# Row 43 - was BKEXToken, now CEJI (swapped with row 42)
$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue "D43" "0.003278"
$ws.Range("E43").Value = "42CEJICEJI"

# Row 44 - LocalTraders
Set-TextValue "D44" "0.008144"

# Row 45 - CoinLion
Set-TextValue "D45" "0.00005445"

# Row 48 - BOLO
Set-TextValue "D48" "0.08918"
